$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old B2 value; its content is re-inserted one row further down (row 3)
$ws.Range("B2").ClearContents()

# Row 3: re-insert the record that used to live in row 2 (A2/B2)
$ws.Range("A3").Value = "EC20230804"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "00218`n"

# Row 4: the extra row introduced by the bug, holding just the B value
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "00218`n"

# Avoid leftover autofit row-height markers from the multi-line text
$ws.Rows(3).AutoFit()
$ws.Rows(4).AutoFit()

# Reflect the final active selection
$ws.Range("G5").Select()
